# fix(publipostage): Correct status name
#
# The shared strings used for the "statut_label" (column B) and
# "statut_name" (column C) values were renamed:
#   bleu                                                   -> noir
#   résultat et / ou publication posté                     -> résultat postés ou publiés
#   pas de résultat ni de publication                      -> pas de résultat postés ni publiés
#   résultat et / ou publication posté dans les 36 mois     -> résultat postés ou publiés dans les 36 mois
#   résultat et / ou publication posté dans les 12 mois     -> résultat postés ou publiés dans les 12 mois
#
# These strings are shared across many rows, so every cell holding one of
# the old values is updated to the corresponding new value (equivalent to a
# find & replace of the shared string across the whole sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu" = "noir";
    "résultat et / ou publication posté" = "résultat postés ou publiés";
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés";
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois";
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois";
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $current = $cell.Value2
        if ($replacements.ContainsKey($current)) {
            $cell.Value = $replacements[$current]
        }
    }
}
